$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new timetracking entry (row 19), matching the formatting of the row above it
# 2022-01-24, 08:00 - 17:00, Tag "Code"
$ws.Range("A18:E18").Copy($ws.Range("A19:E19"))

$ws.Range("A19").Value = 44585
$ws.Range("B19").Value = 0.33333333333333331
$ws.Range("C19").Value = 0.70833333333333337
$ws.Range("D19").Formula = "=C19-B19"
$ws.Range("E19").Value = "Code"

$ws.Range("A20").Select() | Out-Null
